# "files added by kavin"
# Adds a new "Details" worksheet (positioned right after "login") that holds
# a small contact-form style table (headers in row 1, one data row in row 2),
# formatted with a Consolas font + wrapped text, and makes it the active sheet.

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item(1)

# Insert the new sheet right after "login" so tab order is login, Details.
$ws = $wb.Worksheets.Add($null, $loginSheet)
$ws.Name = "Details"

# --- column widths / row heights -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 20.73
$ws.Columns.Item(4).ColumnWidth = 12.73
$ws.Columns.Item(6).ColumnWidth = 12.36
$ws.Columns.Item(7).ColumnWidth = 11.45
$ws.Columns.Item(8).ColumnWidth = 10.82
$ws.Columns.Item(9).ColumnWidth = 13.91

$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(2).RowHeight = 40.5

# --- base formatting for the whole used range -------------------------------------
$dataRange = $ws.Range("A1:I2")
$dataRange.Font.Name = "Consolas"
$dataRange.Font.Size = 10.5
$dataRange.WrapText = $true

# --- header row ---------------------------------------------------------------
$headers = @("fullname", "email", "comName", "mobile", "country", "empCount", "jobTitle", "comments", "exp")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- sample data row ------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "Demouser"
$ws.Cells.Item(2, 2).Value = "demotest@gmail.com"
$ws.Cells.Item(2, 3).Value = "demo"
$ws.Cells.Item(2, 4).Value = 6385667530
$ws.Cells.Item(2, 5).Value = "India"

$expDate = Get-Date -Year 2024 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(2, 6).Value = $expDate
$ws.Cells.Item(2, 6).NumberFormat = "d-mmm"

$ws.Cells.Item(2, 7).Value = "Testing"
$ws.Cells.Item(2, 8).Value = "Demo Test"
$ws.Cells.Item(2, 9).Value = "Thank you."

# --- selection / activation -----------------------------------------------------
[void]$ws.Range("F4").Select()
$ws.Activate()
